# Refresh crypto price/volume snapshot values (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never let Excel coerce numeric-looking
# strings such as "1.00" or "0.0674" into real numbers - the source data is text).
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$updates = @(
    @('D2', '68.097.99'),
    @('E2', '  +0.63%  '),
    @('D3', '3.308.82'),
    @('E3', '  +1.21%  '),
    @('E4', '  +0.00%  '),
    @('D5', '187.02'),
    @('E5', '  +1.62%  '),
    @('D6', '583.49'),
    @('E6', '  +0.47%  '),
    @('E7', '  -0.01%  '),
    @('D8', '0.598'),
    @('E8', '  -0.82%  '),
    @('E9', '  +0.23%  '),
    @('D10', '6.69'),
    @('E10', '  +1.93%  '),
    @('E11', '  +0.37%  '),
    @('D12', '3.884.74'),
    @('E12', '  +1.23%  '),
    @('E13', '  -2.20%  '),
    @('E14', '  +0.89%  '),
    @('D15', '68.215.94'),
    @('E15', '  +0.69%  '),
    @('E16', '  +0.03%  '),
    @('D17', '3.336.93'),
    @('E17', '  +2.06%  '),
    @('D18', '447.85'),
    @('E18', '  +11.64%  '),
    @('D19', '5.72'),
    @('E19', '  +0.07%  '),
    @('D20', '13.54'),
    @('E20', '  +0.80%  '),
    @('E21', '  +2.99%  '),
    @('D22', '74.72'),
    @('E22', '  +5.36%  '),
    @('D23', '0.999'),
    @('E23', '  +0.09%  '),
    @('D24', '3.457.49'),
    @('E24', '  +1.44%  '),
    @('E25', '  +1.37%  '),
    @('E26', '  +1.19%  '),
    @('E27', '  +0.92%  '),
    @('D28', '9.14'),
    @('E28', '  -3.89%  '),
    @('D29', '1.00'),
    @('E29', '  -0.18%  '),
    @('E30', '  +1.54%  '),
    @('D31', '22.89'),
    @('E31', '  +1.22%  '),
    @('D32', '5.37'),
    @('E32', '  -1.89%  '),
    @('E33', '  +0.27%  '),
    @('E35', '  -1.90%  '),
    @('E36', '  +5.55%  '),
    @('D37', '163.76'),
    @('E37', '  +0.24%  '),
    @('D38', '1.85'),
    @('E38', '  -2.44%  '),
    @('E39', '  +0.40%  '),
    @('E40', '  +0.24%  '),
    @('D41', '0.784'),
    @('E41', '  -2.73%  '),
    @('D42', '2.735.30'),
    @('E42', '  +2.65%  '),
    @('D43', '6.41'),
    @('E43', '  +1.93%  '),
    @('D44', '40.67'),
    @('E44', '  -0.02%  '),
    @('D45', '0.0674'),
    @('E45', '  -0.40%  '),
    @('B46', 'InjectiveProtocol'),
    @('C46', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @('D46', '24.79'),
    @('E46', '  +0.96%  '),
    @('B47', 'dogwifhat'),
    @('C47', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'),
    @('D47', '2.41'),
    @('E47', '  -0.72%  '),
    @('D48', '328.55'),
    @('E48', '  -2.46%  '),
    @('D49', '0.0276'),
    @('E49', '  +0.79%  '),
    @('D50', '31.83'),
    @('E50', '  +3.97%  '),
    @('E51', '  +2.63%  ')
)

foreach ($u in $updates) {
    Set-TextValue $u[0] $u[1]
}
